# Itération #4 final - fill in the journal entries for "Iteration #4" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #4")
# "Iteration #3" already has fully filled-in rows; reuse its date-cell
# formatting (borderId + date number format) for the new rows here instead
# of letting auto-detection invent a brand new number format.
$fmtSrc = $wb.Worksheets.Item("Iteration #3")

$fmtSrc.Range("A15").Copy()
$ws.Range("A15:A26").PasteSpecial(-4122)  # xlPasteFormats

function Set-LogRow($Sheet, $Row, $Year, $Month, $Day, $Text, $Hours) {
    $d = Get-Date -Year $Year -Month $Month -Day $Day -Hour 0 -Minute 0 -Second 0
    $Sheet.Cells.Item($Row, 1).Value = $d
    $Sheet.Cells.Item($Row, 2).Value = $Text
    $Sheet.Cells.Item($Row, 3).Value = $Hours
}

Set-LogRow $ws 14 2016 4 10 "ajout de son boom de mort et création de d'icon" 3
Set-LogRow $ws 15 2016 4 11 "ajout d'image de background modification" 2
Set-LogRow $ws 16 2016 4 14 "ajout de la view du marché" 2
Set-LogRow $ws 17 2016 4 17 "travail de création de monstre" 3
Set-LogRow $ws 18 2016 4 18 "création d'un deuxieme monstre" 2
Set-LogRow $ws 19 2016 4 22 "ajout de la nouvelle monnaie pour acheter des chose dans le market avec fonction de calcule" 3
Set-LogRow $ws 20 2016 4 23 "changement de fonctionnement de variable dans un object" 2
Set-LogRow $ws 21 2016 4 29 "modification mineur et optimisation des fonction de calcule pour stopper les erreurs" 3
Set-LogRow $ws 22 2016 4 30 "ajout des bouton du market avec prix et % d'augmentation" 2
Set-LogRow $ws 23 2016 5 1  "tentative de recupérer les donner perdu en changement de view" 3
Set-LogRow $ws 24 2016 5 7  "power points" 4
Set-LogRow $ws 25 2016 5 8  "présentation oral" 3
Set-LogRow $ws 26 2016 5 9  "présentation oral" 2

# Row 19 wraps onto two lines in the real workbook, so the row is a bit taller
# and the description cell keeps its text pinned to the top of the cell.
$ws.Range("B19").VerticalAlignment = -4160  # xlTop
$ws.Rows.Item(19).RowHeight = 16.5

# Effort recap + final comments for this iteration
$ws.Cells.Item(40, 2).Value = 6
$ws.Cells.Item(42, 2).Value = "c'est malheureusement pas ma meilleur itération malgré beaucoup d'heure je semblais ne pas avancer le travail de modification des graphisme prenais beaucoup de temps."

# "Iteration #4" is now the tab shown/selected when the workbook is reopened
$ws.Activate()
$ws.Range("B50").Select()
